$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new monthly rows: 01-07-2021 (row 152) and 01-08-2021 (row 153).
# Column A holds a date-like label ("dd-mm-yyyy" text) that must stay a
# plain text string (matching every other "Serie" cell in column A), not
# get auto-converted into a date serial number by Excel's input parser.
# Trick: write it as a formula producing that literal text, then
# Copy / PasteSpecial-values over itself to collapse it down to a plain
# value in place - this avoids Excel ever treating the original keystrokes
# as date input, so no extra NumberFormat/style ends up applied to the
# cell or left behind unused in the workbook's style table.

$ws.Range("A152").Formula = "=""01-07-2021"""
$ws.Range("A152").Copy()
$ws.Range("A152").PasteSpecial(-4163)

$ws.Range("B152").Value = 5390
$ws.Range("C152").Value = 1151
$ws.Range("D152").Value = 545
$ws.Range("E152").Value = 294
$ws.Range("F152").Value = 1044
$ws.Range("G152").Value = 313
$ws.Range("H152").Value = 1080
$ws.Range("I152").Value = 525
$ws.Range("J152").Value = 437

$ws.Range("A153").Formula = "=""01-08-2021"""
$ws.Range("A153").Copy()
$ws.Range("A153").PasteSpecial(-4163)

$ws.Range("B153").Value = 8149
$ws.Range("C153").Value = 1007
$ws.Range("D153").Value = 924
$ws.Range("E153").Value = 1100
$ws.Range("F153").Value = 1469
$ws.Range("G153").Value = 753
$ws.Range("H153").Value = 1596
$ws.Range("I153").Value = 763
$ws.Range("J153").Value = 538
